# Updated Generate PLC Code
#
# Adds a new "TextListEntry" worksheet (placed right after the existing
# "TextList" sheet) containing the generated PLC alarm text list entries
# for CM_Item_MTR (From/To codes 2001-2006 and their English texts).

$wb = $excel.ActiveWorkbook

# Remember which sheet was active so we can restore the original selection
# after adding the new sheet (Worksheets.Add activates the new sheet).
$originalActiveSheetName = $wb.ActiveSheet.Name

# Add the new worksheet right after the last existing sheet ("TextList"),
# matching the sheet order produced by the PLC code generator.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "TextListEntry"

# ---- Header row (same look as the bold/bordered/centered header used on
# the TextList sheet) ----
$ws.Range("A1").Value = "Parent"
$ws.Range("B1").Value = "From"
$ws.Range("C1").Value = "To"
$ws.Range("D1").Value = "Text [en-US]"

$headerRange = $ws.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# ---- Data rows: Parent, From, To, Text [en-US] ----
# From/To are textual alarm codes (not numbers), so they're entered with a
# leading apostrophe to force text storage, then the cell formatting is
# reset back to the default "Normal" style (clearing the transient
# quote-prefix flag) so no visible/style change remains on the cell.
$data = @(
    @("CM_Item_MTR", "2001", "2001", "Agitator Start"),
    @("CM_Item_MTR", "2002", "2002", "Circulation pump"),
    @("CM_Item_MTR", "2003", "2003", "Antifoam Pump"),
    @("CM_Item_MTR", "2004", "2004", "Base Pump"),
    @("CM_Item_MTR", "2005", "2005", "Feed Pump"),
    @("CM_Item_MTR", "2006", "2006", "Innoculum Pump")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]

    $ws.Cells.Item($row, 2).Value = "'" + $entry[1]
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "'" + $entry[2]
    $ws.Cells.Item($row, 3).Style = "Normal"

    $ws.Cells.Item($row, 4).Value = $entry[3]

    $row++
}

# Restore the originally active sheet/selection.
$wb.Worksheets.Item($originalActiveSheetName).Activate()
